# repull data, push all data, mean calculation
# Update the dSF (column F) values for a handful of rows to reflect
# re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -6
$ws.Range("F5").Value  = -9
$ws.Range("F6").Value  = -3
$ws.Range("F8").Value  = -3
$ws.Range("F9").Value  = 3
$ws.Range("F14").Value = -3
